$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header + country/count data for rows 1-82 (A: country/region name, B: count)
$ws.Cells.Item(1, 1).Value = "Countries/Regions"
$ws.Cells.Item(1, 2).Value = "Count"
$ws.Cells.Item(2, 1).Value = "PEOPLES R CHINA"
$ws.Cells.Item(2, 2).Value = 1920
$ws.Cells.Item(3, 1).Value = "USA"
$ws.Cells.Item(3, 2).Value = 570
$ws.Cells.Item(4, 1).Value = "INDIA"
$ws.Cells.Item(4, 2).Value = 246
$ws.Cells.Item(5, 1).Value = "GERMANY"
$ws.Cells.Item(5, 2).Value = 242
$ws.Cells.Item(6, 1).Value = "SOUTH KOREA"
$ws.Cells.Item(6, 2).Value = 149
$ws.Cells.Item(7, 1).Value = "IRAN"
$ws.Cells.Item(7, 2).Value = 111
$ws.Cells.Item(8, 1).Value = "FRANCE"
$ws.Cells.Item(8, 2).Value = 101
$ws.Cells.Item(9, 1).Value = "JAPAN"
$ws.Cells.Item(9, 2).Value = 100
$ws.Cells.Item(10, 1).Value = "AUSTRALIA"
$ws.Cells.Item(10, 2).Value = 96
$ws.Cells.Item(11, 1).Value = "ENGLAND"
$ws.Cells.Item(11, 2).Value = 82
$ws.Cells.Item(12, 1).Value = "NETHERLANDS"
$ws.Cells.Item(12, 2).Value = 78
$ws.Cells.Item(13, 1).Value = "CANADA"
$ws.Cells.Item(13, 2).Value = 77
$ws.Cells.Item(14, 1).Value = "ITALY"
$ws.Cells.Item(14, 2).Value = 73
$ws.Cells.Item(15, 1).Value = "RUSSIA"
$ws.Cells.Item(15, 2).Value = 72
$ws.Cells.Item(16, 1).Value = "SPAIN"
$ws.Cells.Item(16, 2).Value = 68
$ws.Cells.Item(17, 1).Value = "PORTUGAL"
$ws.Cells.Item(17, 2).Value = 39
$ws.Cells.Item(18, 1).Value = "BELGIUM"
$ws.Cells.Item(18, 2).Value = 38
$ws.Cells.Item(19, 1).Value = "SINGAPORE"
$ws.Cells.Item(19, 2).Value = 38
$ws.Cells.Item(20, 1).Value = "SAUDI ARABIA"
$ws.Cells.Item(20, 2).Value = 36
$ws.Cells.Item(21, 1).Value = "TAIWAN"
$ws.Cells.Item(21, 2).Value = 36
$ws.Cells.Item(22, 1).Value = "SWEDEN"
$ws.Cells.Item(22, 2).Value = 32
$ws.Cells.Item(23, 1).Value = "SWITZERLAND"
$ws.Cells.Item(23, 2).Value = 32
$ws.Cells.Item(24, 1).Value = "POLAND"
$ws.Cells.Item(24, 2).Value = 26
$ws.Cells.Item(25, 1).Value = "MALAYSIA"
$ws.Cells.Item(25, 2).Value = 24
$ws.Cells.Item(26, 1).Value = "THAILAND"
$ws.Cells.Item(26, 2).Value = 23
$ws.Cells.Item(27, 1).Value = "DENMARK"
$ws.Cells.Item(27, 2).Value = 22
$ws.Cells.Item(28, 1).Value = "BRAZIL"
$ws.Cells.Item(28, 2).Value = 20
$ws.Cells.Item(29, 1).Value = "ISRAEL"
$ws.Cells.Item(29, 2).Value = 19
$ws.Cells.Item(30, 1).Value = "FINLAND"
$ws.Cells.Item(30, 2).Value = 18
$ws.Cells.Item(31, 1).Value = "EGYPT"
$ws.Cells.Item(31, 2).Value = 17
$ws.Cells.Item(32, 1).Value = "PAKISTAN"
$ws.Cells.Item(32, 2).Value = 17
$ws.Cells.Item(33, 1).Value = "ROMANIA"
$ws.Cells.Item(33, 2).Value = 17
$ws.Cells.Item(34, 1).Value = "AUSTRIA"
$ws.Cells.Item(34, 2).Value = 14
$ws.Cells.Item(35, 1).Value = "CZECH REPUBLIC"
$ws.Cells.Item(35, 2).Value = 14
$ws.Cells.Item(36, 1).Value = "NIGERIA"
$ws.Cells.Item(36, 2).Value = 14
$ws.Cells.Item(37, 1).Value = "TURKEY"
$ws.Cells.Item(37, 2).Value = 12
$ws.Cells.Item(38, 1).Value = "TURKIYE"
$ws.Cells.Item(38, 2).Value = 12
$ws.Cells.Item(39, 1).Value = "VIETNAM"
$ws.Cells.Item(39, 2).Value = 11
$ws.Cells.Item(40, 1).Value = "CHILE"
$ws.Cells.Item(40, 2).Value = 10
$ws.Cells.Item(41, 1).Value = "SCOTLAND"
$ws.Cells.Item(41, 2).Value = 10
$ws.Cells.Item(42, 1).Value = "IRELAND"
$ws.Cells.Item(42, 2).Value = 9
$ws.Cells.Item(43, 1).Value = "NORWAY"
$ws.Cells.Item(43, 2).Value = 9
$ws.Cells.Item(44, 1).Value = "SOUTH AFRICA"
$ws.Cells.Item(44, 2).Value = 9
$ws.Cells.Item(45, 1).Value = "UKRAINE"
$ws.Cells.Item(45, 2).Value = 9
$ws.Cells.Item(46, 1).Value = "QATAR"
$ws.Cells.Item(46, 2).Value = 8
$ws.Cells.Item(47, 1).Value = "U ARAB EMIRATES"
$ws.Cells.Item(47, 2).Value = 7
$ws.Cells.Item(48, 1).Value = "GREECE"
$ws.Cells.Item(48, 2).Value = 6
$ws.Cells.Item(49, 1).Value = "HUNGARY"
$ws.Cells.Item(49, 2).Value = 6
$ws.Cells.Item(50, 1).Value = "IRAQ"
$ws.Cells.Item(50, 2).Value = 6
$ws.Cells.Item(51, 1).Value = "MEXICO"
$ws.Cells.Item(51, 2).Value = 6
$ws.Cells.Item(52, 1).Value = "SLOVENIA"
$ws.Cells.Item(52, 2).Value = 6
$ws.Cells.Item(53, 1).Value = "LITHUANIA"
$ws.Cells.Item(53, 2).Value = 5
$ws.Cells.Item(54, 1).Value = "NEW ZEALAND"
$ws.Cells.Item(54, 2).Value = 5
$ws.Cells.Item(55, 1).Value = "UZBEKISTAN"
$ws.Cells.Item(55, 2).Value = 5
$ws.Cells.Item(56, 1).Value = "WALES"
$ws.Cells.Item(56, 2).Value = 5
$ws.Cells.Item(57, 1).Value = "BANGLADESH"
$ws.Cells.Item(57, 2).Value = 4
$ws.Cells.Item(58, 1).Value = "BELARUS"
$ws.Cells.Item(58, 2).Value = 4
$ws.Cells.Item(59, 1).Value = "CROATIA"
$ws.Cells.Item(59, 2).Value = 4
$ws.Cells.Item(60, 1).Value = "KAZAKHSTAN"
$ws.Cells.Item(60, 2).Value = 4
$ws.Cells.Item(61, 1).Value = "SLOVAKIA"
$ws.Cells.Item(61, 2).Value = 4
$ws.Cells.Item(62, 1).Value = "ARGENTINA"
$ws.Cells.Item(62, 2).Value = 3
$ws.Cells.Item(63, 1).Value = "ARMENIA"
$ws.Cells.Item(63, 2).Value = 3
$ws.Cells.Item(64, 1).Value = "ESTONIA"
$ws.Cells.Item(64, 2).Value = 3
$ws.Cells.Item(65, 1).Value = "LUXEMBOURG"
$ws.Cells.Item(65, 2).Value = 3
$ws.Cells.Item(66, 1).Value = "AZERBAIJAN"
$ws.Cells.Item(66, 2).Value = 2
$ws.Cells.Item(67, 1).Value = "ETHIOPIA"
$ws.Cells.Item(67, 2).Value = 2
$ws.Cells.Item(68, 1).Value = "INDONESIA"
$ws.Cells.Item(68, 2).Value = 2
$ws.Cells.Item(69, 1).Value = "KUWAIT"
$ws.Cells.Item(69, 2).Value = 2
$ws.Cells.Item(70, 1).Value = "MOROCCO"
$ws.Cells.Item(70, 2).Value = 2
$ws.Cells.Item(71, 1).Value = "SERBIA"
$ws.Cells.Item(71, 2).Value = 2
$ws.Cells.Item(72, 1).Value = "BOSNIA HERCEG"
$ws.Cells.Item(72, 2).Value = 1
$ws.Cells.Item(73, 1).Value = "BOTSWANA"
$ws.Cells.Item(73, 2).Value = 1
$ws.Cells.Item(74, 1).Value = "BULGARIA"
$ws.Cells.Item(74, 2).Value = 1
$ws.Cells.Item(75, 1).Value = "CYPRUS"
$ws.Cells.Item(75, 2).Value = 1
$ws.Cells.Item(76, 1).Value = "GEORGIA"
$ws.Cells.Item(76, 2).Value = 1
$ws.Cells.Item(77, 1).Value = "KENYA"
$ws.Cells.Item(77, 2).Value = 1
$ws.Cells.Item(78, 1).Value = "MAURITIUS"
$ws.Cells.Item(78, 2).Value = 1
$ws.Cells.Item(79, 1).Value = "PERU"
$ws.Cells.Item(79, 2).Value = 1
$ws.Cells.Item(80, 1).Value = "PHILIPPINES"
$ws.Cells.Item(80, 2).Value = 1
$ws.Cells.Item(81, 1).Value = "SRI LANKA"
$ws.Cells.Item(81, 2).Value = 1
$ws.Cells.Item(82, 1).Value = "URUGUAY"
$ws.Cells.Item(82, 2).Value = 1
